$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manualStatus column (I2:I5) from numeric 128 to text "[128]"
$ws.Range("I2").Value = "[128]"
$ws.Range("I3").Value = "[128]"
$ws.Range("I4").Value = "[128]"
$ws.Range("I5").Value = "[128]"

# Move the active selection to I5
$ws.Range("I5").Select()
